$d = $word.ActiveDocument

# 1. Merge the "Andrew Yong " / "Hao" / " Chen " runs (team-member name) into
#    a single run — also clears the spell-check proofErr markers that
#    bracketed "Hao".
$d.Content.Find.Execute("Andrew Yong Hao Chen ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Andrew Yong Hao Chen ", 1)

# 2. Merge the "Game Title: " / "Teleportals" runs into a single run — also
#    clears the spell-check proofErr markers around "Teleportals".
$d.Content.Find.Execute("Game Title: Teleportals", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Game Title: Teleportals", 1)

# 3. Merge the "rather than learning ... more easier to make functions we
#    wanted." runs into a single run — also clears the grammar-check
#    proofErr markers that bracketed "more easier".
$d.Content.Find.Execute("rather than learning 2 engines at the same timeframe. It’s also programming based, so it’s more easier to make functions we wanted.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "rather than learning 2 engines at the same timeframe. It’s also programming based, so it’s more easier to make functions we wanted.", 1)

# 4. Merge the "Cubes to press and hold the " / "Buuttons" / "," runs into a
#    single run — also clears the spell-check proofErr markers around
#    "Buuttons".
$d.Content.Find.Execute("Cubes to press and hold the Buuttons,", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Cubes to press and hold the Buuttons,", 1)

# 5. Merge the "csv " / "maploader" / " / social feature / portal
#    projection" runs into a single run — also clears the spell-check
#    proofErr markers around "maploader".
$d.Content.Find.Execute("csv maploader / social feature / portal projection", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "csv maploader / social feature / portal projection", 1)

# 6. Update my (Isaac's) task list with the new items.
$d.Content.Find.Execute("Isaac: physics / Traps / portal collision / button and door", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Isaac: physics / Traps / portal collision / button and door / audio / coins / level transition", 1)
